# Sincronização de dados: insere um novo registro de avaliação de garantia
# logo antes da antiga linha 22, empurrando as linhas 22..32 para 23..33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insere uma nova linha na posição 22 (as linhas existentes descem uma posição)
$ws.Rows.Item(22).Insert()

# Preenche os dados do novo registro sincronizado
$ws.Cells.Item(22, 1).Value = 5
$ws.Cells.Item(22, 2).Value = ""
$ws.Cells.Item(22, 3).Value = 46006.74939984953
$ws.Cells.Item(22, 4).Value = "ZmYzMjBkNWUtOWQ5YS00MDFiLTk2NmItZTBlMGFlNjE3YmZiOjU3MDE2"
